# Update the CircadiPy cosinor analysis results (re-ran simulations produced
# slightly different fitted values for the two cosinor rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("E2").Value = 25.25000000000051
$ws.Range("H2").Value = [double]"2.64338815386942e-16"
$ws.Range("I2").Value = 0.526277840635985
$ws.Range("K2").Value = 48.32156918624555
$ws.Range("L2").Value = "[37.31897161199863, 59.32416676049248]"
$ws.Range("M2").Value = [double]"4.440892098500626e-16"
$ws.Range("N2").Value = [double]"8.881784197001252e-16"
$ws.Range("O2").Value = 1.314500229429963
$ws.Range("P2").Value = "[1.0755001877154244, 1.5535002711445012]"
$ws.Range("S2").Value = 64.0990291786805
$ws.Range("T2").Value = "[57.82913203170621, 70.36892632565478]"
$ws.Range("W2").Value = 19.96746746746787
$ws.Range("X2").Value = 19.0070070070074
$ws.Range("Y2").Value = 20.92792792792835

# --- Row 3 ---
$ws.Range("E3").Value = 25.45000000000054
$ws.Range("H3").Value = [double]"2.64338815386942e-16"
$ws.Range("I3").ClearContents()
$ws.Range("K3").Value = 48.26824845307563
$ws.Range("L3").Value = "[36.568153406615295, 59.96834349953596]"
$ws.Range("M3").Value = [double]"9.325873406851315e-15"
$ws.Range("N3").Value = [double]"9.325873406851315e-15"
$ws.Range("O3").Value = 1.125815985971117
$ws.Range("P3").Value = "[0.8868159442565782, 1.3648160276856567]"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 66.51369112651726
$ws.Range("T3").Value = "[60.25770679486969, 72.76967545816483]"
$ws.Range("W3").Value = 20.88988988989033
$ws.Range("X3").Value = 19.92182182182224
$ws.Range("Y3").Value = 21.85795795795842
